$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue $ws 'D2' '307.83'
Set-TextValue $ws 'E2' '1.65%'

Set-TextValue $ws 'D3' '37.86'
Set-TextValue $ws 'E3' '6.24%'

Set-TextValue $ws 'D4' '5.077'
Set-TextValue $ws 'E4' '0.82%'

Set-TextValue $ws 'D5' '0.08150'
Set-TextValue $ws 'E5' '3.29%'

Set-TextValue $ws 'D6' '1.980'
Set-TextValue $ws 'E6' '6.85%'

$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue $ws 'D7' '7.897'
Set-TextValue $ws 'E7' '1.54%'

$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D8' '0.9289'
Set-TextValue $ws 'E8' '0.69%'

$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D9' '0.1411'
Set-TextValue $ws 'E9' '2.54%'

$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D10' '0.1947'
Set-TextValue $ws 'E10' '2.59%'

$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D11' '0.09263'
Set-TextValue $ws 'E11' '1.44%'

$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D12' '0.03510'
Set-TextValue $ws 'E12' '1.12%'

$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D13' '0.09886'
Set-TextValue $ws 'E13' '0.49%'

$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D14' '0.001407'
Set-TextValue $ws 'E14' '-0.37%'

$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D15' '0.006239'
Set-TextValue $ws 'E15' '0.30%'

$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D16' '3.944'
Set-TextValue $ws 'E16' '5.63%'

$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws 'D17' '4.174'
Set-TextValue $ws 'E17' '1.49%'

Set-TextValue $ws 'D18' '3.411'
Set-TextValue $ws 'E18' '2.15%'

Set-TextValue $ws 'D19' '0.3454'
Set-TextValue $ws 'E19' '0.38%'

Set-TextValue $ws 'E20' '-4.35%'

Set-TextValue $ws 'D21' '4.812'
Set-TextValue $ws 'E21' '-6.75%'

Set-TextValue $ws 'E22' '18.83%'

Set-TextValue $ws 'D23' '0.04490'
Set-TextValue $ws 'E23' '1.93%'

Set-TextValue $ws 'E24' '0.73%'

Set-TextValue $ws 'E25' '-9.71%'

Set-TextValue $ws 'D27' '0.0001301'
Set-TextValue $ws 'E27' '-0.05%'

Set-TextValue $ws 'D39' '0.02116'
Set-TextValue $ws 'E39' '9.56%'

Set-TextValue $ws 'D40' '0.05150'
Set-TextValue $ws 'E40' '1.45%'

Set-TextValue $ws 'D41' '0.007471'
Set-TextValue $ws 'E41' '-1.42%'

Set-TextValue $ws 'D42' '0.01012'
Set-TextValue $ws 'E42' '-0.45%'

Set-TextValue $ws 'D43' '0.1368'
Set-TextValue $ws 'E43' '1.83%'

Set-TextValue $ws 'D44' '0.002132'
Set-TextValue $ws 'E44' '-1.44%'

Set-TextValue $ws 'D45' '0.009690'
Set-TextValue $ws 'E45' '-0.33%'

Set-TextValue $ws 'D46' '0.00006390'
Set-TextValue $ws 'E46' '3.51%'

Set-TextValue $ws 'D47' '0.00000000751'
Set-TextValue $ws 'E47' '-0.05%'

Set-TextValue $ws 'E48' '-0.64%'

Set-TextValue $ws 'D49' '0.001602'
Set-TextValue $ws 'E49' '-3.56%'

Set-TextValue $ws 'D50' '0.00002102'
Set-TextValue $ws 'E50' '-0.05%'

Set-TextValue $ws 'D51' '0.0002002'
Set-TextValue $ws 'E51' '-0.05%'
